$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Topics")

$ws.Range("B1:C1").EntireColumn.Insert()

$ws.Range("B1").Value = "Source"
$ws.Range("C1").Value = "Type"

$sources = @("Temple CIS","Temple CIS","Temple CIS","Temple CIS","Temple CIS","Temple CIS","Temple CIS","Temple CIS","Temple CIS","Temple CIS","Temple CIS")
$types   = @("CS Core","CS Core","CS Core","CS Core","CS Core","CS Core","CS Electives","CS Core","CS Core","CS Core","CS Core")

for ($i = 0; $i -lt 11; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $sources[$i]
    $ws.Range("C$row").Value = $types[$i]
}

$ws.Range("A13").Value = 1012
$ws.Range("B13").Value = "Temple CIS"
$ws.Range("C13").Value = "KA Core"
$ws.Range("D13").Value = "Machine Learning"
$ws.Range("E13").Value = 4000
$ws.Range("F13").Value = "Logistics Regression"

$ws.Range("H2").Value = 1052
$ws.Range("H3").Value = 1053
$ws.Range("H4").Value = 1054
$ws.Range("H5").Value = 1055
$ws.Range("H6").Value = 1056
$ws.Range("H7").Value = 1057
$ws.Range("H8").Value = 1058
$ws.Range("H9").Value = 1059
$ws.Range("H10").Value = 1060
$ws.Range("H11").Value = 1061
$ws.Range("H12").Value = 1062

Write-Output "done"
